$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 245, shifting existing rows 245:303 down to 246:304
$ws.Rows(245).EntireRow.Insert()

# Populate the newly inserted row 245 with the new data record
$ws.Range("A245").Value = 10
$ws.Range("B245").Value = "Vega Modelo de Temuco"
$ws.Range("C245").Value = "La Araucanía"
$ws.Range("D245").NumberFormat = $ws.Range("D246").NumberFormat
$ws.Range("D245").Value = 45211
$ws.Range("E245").Value = 9
$ws.Range("F245").Value = "Fruta"
$ws.Range("G245").Value = 100104
$ws.Range("H245").Value = "Frutos de pepita"
$ws.Range("I245").Value = 100104001
$ws.Range("J245").Value = "Granada"
$ws.Range("K245").Value = "Wonderfull"
$ws.Range("L245").Value = "Primera"
$ws.Range("M245").Value = 270
$ws.Range("N245").Value = 17000
$ws.Range("O245").Value = 17000
$ws.Range("P245").Value = 17000
$ws.Range("Q245").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R245").Value = "Provincia de Limarí"
$ws.Range("S245").Value = 1700
$ws.Range("T245").Value = 10
